$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Sanguinaria canadensis",
    "Aquilegia canadensis",
    "Caltha palustris",
    "Caltha palustris",
    "Dicentra cucullaria",
    "Asarum canadense",
    "Hepatica americana",
    "Hepatica americana",
    "Arisaema triphyllum",
    "Podophyllum peltatum",
    "Phlox divaricata",
    "Phlox divaricata",
    "Claytonia Virginica",
    "Trillium grandiflorum",
    "Trillium grandiflorum",
    "Erythronium americanum",
    "Erythronium americanum",
    "Erythronium americanum",
    "Anemone blanda",
    "Anemone blanda",
    "Monarda didyma",
    "Monarda didyma",
    "Rudbeckia hirta",
    "Ranunculus",
    "Ranunculus",
    "Asclepias tuberosa",
    "Potentilla",
    "Oenothera",
    "Gentiana",
    "Gentiana",
    "Polemonium caeruleum",
    "Polemonium caeruleum",
    "Eschscholzia californica",
    "Dodecatheon",
    "Cimicifuga",
    "Lobelia cardinalis"
)

$startRow = 4
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("A8").Select()
